$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "two_by_two_PriceinIntput"

$ws.Range("B1").Value = "'benchmark"
$ws.Range("C1").Value = "'RA=157"
$ws.Range("D1").Value = "'PX=1"
$ws.Range("E1").Value = "'PL=1"
$ws.Range("F1").Value = "'Itax=0.1"
$ws.Range("G1").Value = "'Otax=0.1"
$ws.Range("A2").Value = "'X"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1.0533624393223964
$ws.Range("D2").Value = 1.0533624393076026
$ws.Range("E2").Value = 1.0533624392983827
$ws.Range("F2").Value = 1.0692230525392665
$ws.Range("G2").Value = 1.0306253077634344
$ws.Range("A3").Value = "'Y"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1.0432700717339818
$ws.Range("D3").Value = 1.0432700717206651
$ws.Range("E3").Value = 1.043270071718168
$ws.Range("F3").Value = 1.0159922067315308
$ws.Range("G3").Value = 1.0873120422296021
$ws.Range("A4").Value = "'U"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.0500286076528291
$ws.Range("D4").Value = 1.0500286076372356
$ws.Range("E4").Value = 1.0500286076315015
$ws.Range("F4").Value = 1.0513947937746249
$ws.Range("G4").Value = 1.0489563308826293
$ws.Range("A5").Value = "'PX"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0.99683505738742828
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1.0442749417809107
$ws.Range("F5").Value = 0.96924826808406228
$ws.Range("G5").Value = 1.080588241014147
$ws.Range("A6").Value = "'PY"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1.006478222756964
$ws.Range("D6").Value = 1.0096737823300179
$ws.Range("E6").Value = 1.0543770302540776
$ws.Range("F6").Value = 1.0200300602347367
$ws.Range("G6").Value = 1.0242520502029886
$ws.Range("A7").Value = "'PU"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1.0031749912758545
$ws.Range("E7").Value = 1.0475905056086889
$ws.Range("F7").Value = 0.98568358717211357
$ws.Range("G7").Value = 1.0617044348776998
$ws.Range("A8").Value = "'PL"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.95457146145108884
$ws.Range("D8").Value = 0.95760221755176367
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("A9").Value = "'PK"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1.0500286076528458
$ws.Range("D9").Value = 1.0533624393068366
$ws.Range("E9").Value = 1.0999999999351249
$ws.Range("F9").Value = 1.0363425908479043
$ws.Range("G9").Value = 1.0440764845854245
$ws.Range("A10").Value = "'SX"
$ws.Range("B10").Value = 110
$ws.Range("C10").Value = 110
$ws.Range("D10").Value = 110
$ws.Range("E10").Value = 110
$ws.Range("F10").Value = 110
$ws.Range("G10").Value = 110
$ws.Range("A11").Value = "'SY"
$ws.Range("B11").Value = 54
$ws.Range("C11").Value = 54
$ws.Range("D11").Value = 54
$ws.Range("E11").Value = 54
$ws.Range("F11").Value = 54
$ws.Range("G11").Value = 54
$ws.Range("A12").Value = "'SU"
$ws.Range("B12").Value = 164
$ws.Range("C12").Value = 164
$ws.Range("D12").Value = 164
$ws.Range("E12").Value = 164
$ws.Range("F12").Value = 164
$ws.Range("G12").Value = 164
$ws.Range("A13").Value = "'DXL"
$ws.Range("B13").Value = 50
$ws.Range("C13").Value = 52.213747091206962
$ws.Range("D13").Value = 52.213747089924226
$ws.Range("E13").Value = 52.213747088526816
$ws.Range("F13").Value = 52.868087335139379
$ws.Range("G13").Value = 53.047059029882369
$ws.Range("A14").Value = "'DXK"
$ws.Range("B14").Value = 50
$ws.Range("C14").Value = 47.467042807627351
$ws.Range("D14").Value = 47.46704280902668
$ws.Range("E14").Value = 47.467042810551135
$ws.Range("F14").Value = 46.762927451330512
$ws.Range("G14").Value = 46.573667250728739
$ws.Range("A15").Value = "'DYL"
$ws.Range("B15").Value = 24
$ws.Range("C15").Value = 25.305048727435071
$ws.Range("D15").Value = 25.305048726675256
$ws.Range("E15").Value = 25.305048725847506
$ws.Range("F15").Value = 24.480721438762444
$ws.Range("G15").Value = 24.582049170514459
$ws.Range("A16").Value = "'DYK"
$ws.Range("B16").Value = 30
$ws.Range("C16").Value = 28.755737188715777
$ws.Range("D16").Value = 28.755737189406524
$ws.Range("E16").Value = 28.755737190159014
$ws.Range("F16").Value = 29.527785568878645
$ws.Range("G16").Value = 29.43037403561884
$ws.Range("A17").Value = "'DUX"
$ws.Range("B17").Value = 110
$ws.Range("C17").Value = 110.3492490405403
$ws.Range("D17").Value = 110.34924904034398
$ws.Range("E17").Value = 110.34924904012642
$ws.Range("F17").Value = 111.86524460134066
$ws.Range("G17").Value = 108.07769639007446
$ws.Range("A18").Value = "'DUY"
$ws.Range("B18").Value = 54
$ws.Range("C18").Value = 53.652427622403614
$ws.Range("D18").Value = 53.652427622598076
$ws.Range("E18").Value = 53.652427622813583
$ws.Range("F18").Value = 52.181710894781091
$ws.Range("G18").Value = 55.974542078075267
$ws.Range("A19").Value = "'RA"
$ws.Range("B19").Value = 164
$ws.Range("C19").Value = 172.20469165071995
$ws.Range("D19").Value = 172.75144004633617
$ws.Range("E19").Value = 180.39999999442588
$ws.Range("F19").Value = 169.96018504007139
$ws.Range("G19").Value = 182.64378039376402
$ws.Range("A20").Value = "'DU"
$ws.Range("B20").Value = 164
$ws.Range("C20").Value = 172.20469165071995
$ws.Range("D20").Value = 172.20469165267772
$ws.Range("E20").Value = 172.20469165058611
$ws.Range("F20").Value = 172.42874615339829
$ws.Range("G20").Value = 172.02883815287367
$ws.Range("A21").Value = "'CWI"
$ws.Range("B21").Value = 1.0933333333333333
$ws.Range("C21").Value = 1.1480312776714663
$ws.Range("D21").Value = 1.1480312776845181
$ws.Range("E21").Value = 1.1480312776705741
$ws.Range("F21").Value = 1.1495249743559885
$ws.Range("G21").Value = 1.1468589210191578
$ws.Range("A22").Value = "'PX/PX"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("A23").Value = "'PY/PX"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 1.009673782335474
$ws.Range("D23").Value = 1.0096737823300179
$ws.Range("E23").Value = 1.0096737823239719
$ws.Range("F23").Value = 1.0523929666142775
$ws.Range("G23").Value = 0.94786525646597253
$ws.Range("A24").Value = "'PU/PX"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 1.0031749912777612
$ws.Range("D24").Value = 1.0031749912758545
$ws.Range("E24").Value = 1.0031749912739683
$ws.Range("F24").Value = 1.0169567691057517
$ws.Range("G24").Value = 0.98252451264995777
$ws.Range("A25").Value = "'PL/PX"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 0.95760221751519581
$ws.Range("D25").Value = 0.95760221755176367
$ws.Range("E25").Value = 0.95760221756790975
$ws.Range("F25").Value = 1.0317274045552081
$ws.Range("G25").Value = 0.92542187860705039
$ws.Range("A26").Value = "'PK/PX"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 1.0533624393235435
$ws.Range("D26").Value = 1.0533624393068366
$ws.Range("E26").Value = 1.0533624392625762
$ws.Range("F26").Value = 1.0692230514855281
$ws.Range("G26").Value = 0.96621122177448859
$ws.Range("A27").Value = "'RA/PX"
$ws.Range("B27").Value = 164
$ws.Range("C27").Value = 172.75144004470056
$ws.Range("D27").Value = 172.75144004633617
$ws.Range("E27").Value = 172.75144004391314
$ws.Range("F27").Value = 175.35258058911575
$ws.Range("G27").Value = 169.02255036789066
